$wb = $excel.ActiveWorkbook

# --- Metadata: bump "Last Updated" timestamp ---
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("A2").Value = "05 Nov 2025, 02:56 PM"

# --- Industry Analysis: refresh "1 Year" (column F) figures ---
$wsInd = $wb.Worksheets.Item("Industry Analysis")
$wsInd.Range("F2").Value = 21.0016
$wsInd.Range("F3").Value = -16.2396
$wsInd.Range("F4").Value = 27.1317
$wsInd.Range("F5").Value = -50.6494
$wsInd.Range("F6").Value = 53.2813
$wsInd.Range("F7").Value = -8.106199999999999
$wsInd.Range("F8").Value = -9.552099999999999
$wsInd.Range("F9").Value = 36.3756
$wsInd.Range("F10").Value = -6.1314
$wsInd.Range("F11").Value = 31.9081
$wsInd.Range("F12").Value = -18.4955
$wsInd.Range("F13").Value = 14.0155
$wsInd.Range("F14").Value = -36.0718
$wsInd.Range("F15").Value = -0.1622
$wsInd.Range("F16").Value = 0.1459
$wsInd.Range("F17").Value = -22.0012
$wsInd.Range("F18").Value = 1.0561
$wsInd.Range("F19").Value = -27.708
$wsInd.Range("F20").Value = 47.7309
$wsInd.Range("F21").Value = 12.0959
$wsInd.Range("F22").Value = 95.1491
$wsInd.Range("F23").Value = -50.2657
$wsInd.Range("F24").Value = -13.3427
$wsInd.Range("F25").Value = -9.9316
$wsInd.Range("F26").Value = 5.8244
$wsInd.Range("F27").Value = -32.7692
$wsInd.Range("F28").Value = -24.8224
$wsInd.Range("F29").Value = -18.4191
$wsInd.Range("F30").Value = 25.8569
$wsInd.Range("F31").Value = 58.4712
$wsInd.Range("F32").Value = -3.3862
$wsInd.Range("F33").Value = -6.3282
$wsInd.Range("F34").Value = 27.7203
$wsInd.Range("F35").Value = 4.4873
$wsInd.Range("F36").Value = -4.9458
$wsInd.Range("F37").Value = 3.6074
$wsInd.Range("F38").Value = -23.3973
$wsInd.Range("F39").Value = 8.7355
$wsInd.Range("F40").Value = -5.8541
$wsInd.Range("F41").Value = -8.3934
$wsInd.Range("F42").Value = 20.3818
$wsInd.Range("F43").Value = 14.3164
$wsInd.Range("F44").Value = -12.6846
$wsInd.Range("F45").Value = 28.4075
$wsInd.Range("F46").Value = -1.1135
$wsInd.Range("F47").Value = -37.1997
$wsInd.Range("F48").Value = -29.8569
$wsInd.Range("F49").Value = -27.5511
$wsInd.Range("F50").Value = -49.7478
$wsInd.Range("F51").Value = -51.8002
$wsInd.Range("F52").Value = -38.5254
$wsInd.Range("F53").Value = -12.4886
$wsInd.Range("F54").Value = -5.0725
$wsInd.Range("F55").Value = -17.7445
$wsInd.Range("F56").Value = -26.636
$wsInd.Range("F57").Value = -29.3361
$wsInd.Range("F58").Value = -11.9574
$wsInd.Range("F59").Value = -24.5687
$wsInd.Range("F60").Value = -12.3
$wsInd.Range("F61").Value = -10.9446
$wsInd.Range("F62").Value = -17.1229
$wsInd.Range("F63").Value = -9.5038
$wsInd.Range("F64").Value = 54.2749
$wsInd.Range("F65").Value = -43.4736
$wsInd.Range("F66").Value = 13.2687
$wsInd.Range("F67").Value = 12.7149
$wsInd.Range("F68").Value = 24.8057
$wsInd.Range("F69").Value = -17.0328
$wsInd.Range("F70").Value = -6.8927
$wsInd.Range("F71").Value = 13.6034
$wsInd.Range("F72").Value = 3.9995
$wsInd.Range("F73").Value = -16.226
$wsInd.Range("F74").Value = -16.2448
$wsInd.Range("F75").Value = 28.6924
$wsInd.Range("F76").Value = 48.9752

# --- Stock List: a new ticker (CAPTRU-RE1) appears at the top, pushing every
#     other row down by one and dropping the previous last row off the list ---
$wsStock = $wb.Worksheets.Item("Stock List")
$wsStock.Rows.Item(2).Insert()
$wsStock.Range("A2:H2").ClearFormats()
$wsStock.Range("A2").Value = "📋"
$wsStock.Range("B2").Value = "CAPTRU-RE1"
$wsStock.Range("C2").Value = "CAPTRU-RE1"
$wsStock.Range("D2").Value = 5.67
$wsStock.Range("E2").Value = -11.9565
$wsStock.Range("F2").Value = "N/A"
$wsStock.Range("G2").Value = "N/A"
$wsStock.Range("H2").Value = 0
$wsStock.Rows.Item(77).Delete()
